$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 3.2
$ws.Range("I2").Value = 2.84
$ws.Range("J2").Value = 3.55
$ws.Range("L2").Value = 1.24
$ws.Range("M2").Value = 1.01
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 1.74
$ws.Range("T2").Value = 1.44
$ws.Range("U2").Value = 1.9
$ws.Range("V2").Value = 1.54
$ws.Range("W2").Value = 1.46
$ws.Range("X2").Value = 26
$ws.Range("Y2").Value = 18
$ws.Range("Z2").Value = 28
$ws.Range("AA2").Value = 55
$ws.Range("AB2").Value = 20
$ws.Range("AC2").Value = 12.5
$ws.Range("AD2").Value = 18
$ws.Range("AE2").Value = 40
$ws.Range("AF2").Value = 30
$ws.Range("AG2").Value = 18.5
$ws.Range("AH2").Value = 23
$ws.Range("AI2").Value = 50
$ws.Range("AJ2").Value = 65
$ws.Range("AK2").Value = 40
$ws.Range("AL2").Value = 50
$ws.Range("AM2").Value = 100
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3 updates
$ws.Range("L3").Value = 1.41
$ws.Range("N3").Value = 3.7
$ws.Range("P3").Value = 1.9
$ws.Range("Q3").Value = 2.04
$ws.Range("R3").Value = 1.34
$ws.Range("S3").Value = 3.55

# Row 4 updates
$ws.Range("F4").Value = 2.26
$ws.Range("I4").Value = 3.55
$ws.Range("K4").Value = 4.2
$ws.Range("P4").Value = 2.18

# Row 6 updates
$ws.Range("P6").Value = 1.97
